$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the hours worked value for Monday of the week starting 2018-02-19 (row 7)
$ws.Range("B7").Value = 6.75

# Move the active selection to F12, matching the saved cursor position
$ws.Range("F12").Select()
